$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.215.50'

$ws.Range("D3").Value = '1.916.16'
$ws.Range("E3").Value = '  -3.70%  '

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = '  -1.01%  '

$ws.Range("D5").Value = "'327.60"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").Value = "'0.4680"
$ws.Range("E7").Value = '  -6.05%  '

$ws.Range("D8").Value = "'0.4014"
$ws.Range("E8").Value = '  -4.67%  '

$ws.Range("D9").Value = "'53.15"
$ws.Range("E9").Value = '  -2.57%  '

$ws.Range("D10").Value = "'0.08398"
$ws.Range("E10").Value = '  -10.00%  '

$ws.Range("D11").Value = "'1.045"
$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("D12").Value = "'22.12"
$ws.Range("E12").Value = '  -4.06%  '

$ws.Range("D13").Value = '1.887.41'
$ws.Range("E13").Value = '  -4.65%  '

$ws.Range("D14").Value = "'7.424"
$ws.Range("E14").Value = '  -7.04%  '

$ws.Range("D15").Value = "'6.073"
$ws.Range("E15").Value = '  -6.03%  '

$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").Value = "'89.71"
$ws.Range("E17").Value = '  -3.28%  '

$ws.Range("D18").Value = "'0.00001061"
$ws.Range("E18").Value = '  -4.70%  '

$ws.Range("D19").Value = "'0.06594"
$ws.Range("E19").Value = '  -2.06%  '

$ws.Range("D20").Value = "'17.97"
$ws.Range("E20").Value = '  -7.77%  '

$ws.Range("E21").Value = '  -0.97%  '

$ws.Range("D22").Value = "'5.726"
$ws.Range("E22").Value = '  -4.40%  '

$ws.Range("D23").Value = '28.191.69'
$ws.Range("E23").Value = '  -3.27%  '

$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = '  -6.15%  '

$ws.Range("D25").Value = "'2.278"
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("D26").Value = '2.141.56'
$ws.Range("E26").Value = '  -3.68%  '

$ws.Range("D27").Value = "'153.38"
$ws.Range("E27").Value = '  -2.13%  '

$ws.Range("D28").Value = "'19.99"
$ws.Range("E28").Value = '  -4.08%  '

$ws.Range("D29").Value = "'2.126"
$ws.Range("E29").Value = '  -6.32%  '

$ws.Range("D30").Value = "'5.730"
$ws.Range("E30").Value = '  -8.96%  '

$ws.Range("D31").Value = "'123.46"
$ws.Range("E31").Value = '  -3.32%  '

$ws.Range("D32").Value = "'0.9751"
$ws.Range("E32").Value = '  -7.34%  '

$ws.Range("D33").Value = "'0.09605"
$ws.Range("E33").Value = '  -2.48%  '

$ws.Range("D34").Value = "'1.450"
$ws.Range("E34").Value = '  -6.08%  '

$ws.Range("D35").Value = "'5.545"
$ws.Range("E35").Value = '  -4.85%  '

$ws.Range("D36").Value = "'3.627"
$ws.Range("E36").Value = '  -2.99%  '

$ws.Range("D37").Value = "'8.815"
$ws.Range("E37").Value = '  -3.39%  '

$ws.Range("D38").Value = "'0.02301"
$ws.Range("E38").Value = '  -5.29%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = "'1.258"
$ws.Range("E39").Value = '  -4.86%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.06172"
$ws.Range("E40").Value = '  -3.74%  '

$ws.Range("D41").Value = "'0.6132"
$ws.Range("E41").Value = '  -5.68%  '

$ws.Range("D42").Value = "'11.02"
$ws.Range("E42").Value = '  -4.23%  '

$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = '  -0.90%  '

$ws.Range("D44").Value = "'0.1901"
$ws.Range("E44").Value = '  -5.08%  '

$ws.Range("D45").Value = "'1.301"
$ws.Range("E45").Value = '  -4.57%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.5827"
$ws.Range("E46").Value = '  -6.54%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'12.73"
$ws.Range("E47").Value = '  -5.56%  '

$ws.Range("D48").Value = "'2.024"
$ws.Range("E48").Value = '  -7.74%  '

$ws.Range("E49").Value = '  -1.53%  '

$ws.Range("D50").Value = "'0.06873"
$ws.Range("E50").Value = '  -1.85%  '

$ws.Range("D51").Value = "'109.96"
$ws.Range("E51").Value = '  -3.17%  '
